$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 419.75
$ws.Range("I2").Value = 387.66666
$ws.Range("K2").Value = 387.66666
$ws.Range("M2").Value = -274.66666

$ws.Range("H51").Value = 2638.4614
$ws.Range("I51").Value = 2500
$ws.Range("K51").Value = 2500
$ws.Range("M51").Value = -2016

$ws.Range("H98").Value = 763.5
$ws.Range("I98").Value = 763.5
$ws.Range("K98").Value = 763.5
$ws.Range("M98").Value = 734.5

$ws.Range("H122").Value = 763.5
$ws.Range("I122").Value = 763.5
$ws.Range("K122").Value = 2290.5
$ws.Range("M122").Value = 159.5

$ws.Range("H135").Value = 825.3077
$ws.Range("I135").Value = 560.75
$ws.Range("K135").Value = 5046.75
$ws.Range("M135").Value = -2511.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 40000
$ws.Range("J43").Value = 40000
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40626

$ws.Range("H45").Value = 14942.75
$ws.Range("I45").Value = 11521.3125
$ws.Range("J45").Value = 28628.5
$ws.Range("K45").Value = 11521.3125
$ws.Range("L45").Value = 28628.5
$ws.Range("M45").Value = -11144.3125
$ws.Range("N45").Value = -29382.5

$ws.Range("H46").Value = 3999.5
$ws.Range("I46").Value = 3999.5
$ws.Range("K46").Value = 3999.5
$ws.Range("M46").Value = -3680.5

$ws.Range("H61").Value = 3367.6843
$ws.Range("I61").Value = 3165.889
$ws.Range("K61").Value = 3165.889
$ws.Range("M61").Value = -2953.889

$ws.Range("H74").Value = 10098.583
$ws.Range("I74").Value = 1962.2858
$ws.Range("J74").Value = 21489.4
$ws.Range("K74").Value = 1962.2858
$ws.Range("L74").Value = 21489.4
$ws.Range("M74").Value = -1088.2858
$ws.Range("N74").Value = -23237.4

$ws.Range("H77").Value = 10098.583
$ws.Range("I77").Value = 1962.2858
$ws.Range("J77").Value = 21489.4
$ws.Range("K77").Value = 9811.429
$ws.Range("L77").Value = 107447
$ws.Range("M77").Value = -5443.429
$ws.Range("N77").Value = -116183

$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680

$ws.Range("H109").Value = 60000
$ws.Range("J109").Value = 60000
$ws.Range("L109").Value = 60000
$ws.Range("N109").Value = -62774

$ws.Range("H110").Value = 7439.231
$ws.Range("I110").Value = 8461.32
$ws.Range("J110").Value = 5614.0713
$ws.Range("K110").Value = 8461.32
$ws.Range("L110").Value = 5614.0713
$ws.Range("M110").Value = -6416.32
$ws.Range("N110").Value = -9704.0713

$ws.Range("H122").Value = 2830.4546
$ws.Range("I122").Value = 2481.5715
$ws.Range("J122").Value = 3441
$ws.Range("K122").Value = 7444.7145
$ws.Range("L122").Value = 10323
$ws.Range("M122").Value = -4994.7145
$ws.Range("N122").Value = -15223

$ws.Range("H136").Value = 3367.6843
$ws.Range("I136").Value = 3165.889
$ws.Range("K136").Value = 9497.667000000001
$ws.Range("M136").Value = -6947.667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1254.1818
$ws.Range("I22").Value = 324.66666
$ws.Range("K22").Value = 324.66666
$ws.Range("M22").Value = 25.33334000000002

$ws.Range("H63").Value = 80000
$ws.Range("J63").Value = 80000
$ws.Range("L63").Value = 80000
$ws.Range("N63").Value = -81372

$ws.Range("H66").Value = 80000
$ws.Range("J66").Value = 80000
$ws.Range("L66").Value = 240000
$ws.Range("N66").Value = -246864

$ws.Range("H99").Value = 3399.1304
$ws.Range("I99").Value = 3287.1052
$ws.Range("K99").Value = 3287.1052
$ws.Range("M99").Value = -1789.1052

$ws.Range("H126").Value = 3399.1304
$ws.Range("I126").Value = 3287.1052
$ws.Range("K126").Value = 9861.3156
$ws.Range("M126").Value = -7391.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4064.1428
$ws.Range("J104").Value = 5257.4
$ws.Range("L104").Value = 15772.2
$ws.Range("N104").Value = -21014.2

$ws.Range("H124").Value = 24640.834
$ws.Range("I124").Value = 7845
$ws.Range("J124").Value = 28000
$ws.Range("K124").Value = 23535
$ws.Range("L124").Value = 84000
$ws.Range("M124").Value = -18625
$ws.Range("N124").Value = -93820

$ws.Range("H129").Value = 1543
$ws.Range("I129").Value = 949.6667
$ws.Range("J129").Value = 1988
$ws.Range("K129").Value = 2849.0001
$ws.Range("L129").Value = 5964
$ws.Range("M129").Value = 2150.9999
$ws.Range("N129").Value = -15964

$ws.Range("H131").Value = 53893.812
$ws.Range("I131").Value = 286594.16
$ws.Range("J131").Value = 8646.527
$ws.Range("K131").Value = 859782.48
$ws.Range("L131").Value = 25939.581
$ws.Range("M131").Value = -854742.48
$ws.Range("N131").Value = -36019.581

$ws.Range("H134").Value = 7704.1665
$ws.Range("I134").Value = 4605.5557
$ws.Range("K134").Value = 13816.6671
$ws.Range("M134").Value = -8746.667099999999

$ws.Range("H137").Value = 5676.278
$ws.Range("I137").Value = 5526.091
$ws.Range("J137").Value = 5912.2856
$ws.Range("K137").Value = 16578.273
$ws.Range("L137").Value = 17736.8568
$ws.Range("M137").Value = -11478.273
$ws.Range("N137").Value = -27936.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24975
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H64").Value = 58950
$ws.Range("J64").Value = 58950
$ws.Range("L64").Value = 58950
$ws.Range("N64").Value = -59446

$ws.Range("H67").Value = 58950
$ws.Range("J67").Value = 58950
$ws.Range("L67").Value = 58950
$ws.Range("N67").Value = -60666

$ws.Range("H126").Value = 18500.38
$ws.Range("J126").Value = 3735.2222
$ws.Range("L126").Value = 11205.6666
$ws.Range("N126").Value = -16145.6666

$ws.Range("H132").Value = 215564.08
$ws.Range("J132").Value = 2357.8
$ws.Range("L132").Value = 7073.400000000001
$ws.Range("N132").Value = -12133.4

$ws.Range("H141").Value = 38000
$ws.Range("I141").Value = 38000
$ws.Range("K141").Value = 38000
$ws.Range("M141").Value = -32820

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9937.223
$ws.Range("I7").Value = 11114.615
$ws.Range("K7").Value = 11114.615
$ws.Range("M7").Value = -11002.615

$ws.Range("H41").Value = 34000
$ws.Range("J41").Value = 34000
$ws.Range("L41").Value = 34000
$ws.Range("N41").Value = -34876

$ws.Range("H45").Value = 19995
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 19995
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 19995
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -20809

$ws.Range("H55").Value = 147.72223
$ws.Range("I55").Value = 145.84616
$ws.Range("J55").Value = 152.6
$ws.Range("K55").Value = 145.84616
$ws.Range("L55").Value = 152.6
$ws.Range("M55").Value = 27.15384
$ws.Range("N55").Value = -498.6

$ws.Range("H100").Value = 67421.375
$ws.Range("I100").Value = 89086.836
$ws.Range("K100").Value = 89086.836
$ws.Range("M100").Value = -88545.836

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H126").Value = 9937.223
$ws.Range("I126").Value = 11114.615
$ws.Range("K126").Value = 33343.845
$ws.Range("M126").Value = -30873.845

$ws.Range("H132").Value = 5186.3687
$ws.Range("I132").Value = 4856.9287
$ws.Range("J132").Value = 6108.8
$ws.Range("K132").Value = 14570.7861
$ws.Range("L132").Value = 18326.4
$ws.Range("M132").Value = -12040.7861
$ws.Range("N132").Value = -23386.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20328.477
$ws.Range("J62").Value = 23632.834
$ws.Range("L62").Value = 23632.834
$ws.Range("N62").Value = -24880.834

$ws.Range("H65").Value = 20328.477
$ws.Range("J65").Value = 23632.834
$ws.Range("L65").Value = 118164.17
$ws.Range("N65").Value = -124404.17

$ws.Range("H132").Value = 2280.926
$ws.Range("I132").Value = 2011.44
$ws.Range("K132").Value = 6034.32
$ws.Range("M132").Value = -3504.32

$ws.Range("H140").Value = 55670.715
$ws.Range("J140").Value = 55670.715
$ws.Range("L140").Value = 55670.715
$ws.Range("N140").Value = -66030.715

$ws.Range("H141").Value = 104815.57
$ws.Range("J141").Value = 104815.57
$ws.Range("L141").Value = 104815.57
$ws.Range("N141").Value = -115175.57
